# aggiornamento fino a 02/05
# Extend the daily COVID data table (Sheet1) with six new rows (239-244),
# covering dates 2021-04-27 .. 2021-05-02 (Excel serials 44313..44318).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of the last existing data row's date cell
# (A238, which carries the bold/centered/bordered date-time number format)
# down across the new date cells A239:A244, so the new rows look like the
# existing ones (style index s="2" in the OOXML).
$ws.Range("A238").Copy($ws.Range("A239:A244"))

# Row 239: 2021-04-27
$ws.Range("A239").Value = 44313
$ws.Range("B239").Value = 0
$ws.Range("C239").Value = 6
$ws.Range("D239").Value = 91.37983551629607

# Row 240: 2021-04-28
$ws.Range("A240").Value = 44314
$ws.Range("B240").Value = 0
$ws.Range("C240").Value = 5
$ws.Range("D240").Value = 76.14986293024673

# Row 241: 2021-04-29
$ws.Range("A241").Value = 44315
$ws.Range("B241").Value = 0
$ws.Range("C241").Value = 5
$ws.Range("D241").Value = 76.14986293024673

# Row 242: 2021-04-30
$ws.Range("A242").Value = 44316
$ws.Range("B242").Value = 2
$ws.Range("C242").Value = 5
$ws.Range("D242").Value = 76.14986293024673

# Row 243: 2021-05-01
$ws.Range("A243").Value = 44317
$ws.Range("B243").Value = 0
$ws.Range("C243").Value = 3
$ws.Range("D243").Value = 45.68991775814803

# Row 244: 2021-05-02
$ws.Range("A244").Value = 44318
$ws.Range("B244").Value = 1
$ws.Range("C244").Value = 4
$ws.Range("D244").Value = 60.91989034419738
